$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 2: Salchipapa -> Salchipapas, Existencias 4 -> 10
$ws.Range("A2").Value = "Salchipapas"
$ws.Range("D2").Value = 10

# Update row 3: Precio 1 -> 2, Existencias 8 -> 10
$ws.Range("B3").Value = 2
$ws.Range("D3").Value = 10

# Update row 4: Gorro de baño -> gorros, Precio 2 -> 3, Existencias 9 -> 10
$ws.Range("A4").Value = "gorros"
$ws.Range("B4").Value = 3
$ws.Range("D4").Value = 10

# Delete row 5 (Chochos con tostado) entirely
$ws.Rows.Item(5).Delete()
